$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added at the top of the historical table
# (most-recent-first ordering), which pushes the existing rows 469..482
# down to 470..483 (dimension grows from T482 to T483).
$ws.Rows(469).Insert()

# Populate the newly inserted row 469 with the new record's data.
$ws.Range("A469").Value = 5
$ws.Range("B469").Value = "Macroferia Regional de Talca"
$ws.Range("C469").Value = "Maule"
$ws.Range("D469").Value = 45239
$ws.Range("E469").Value = 7
$ws.Range("F469").Value = "Fruta"
$ws.Range("G469").Value = 100108
$ws.Range("H469").Value = "Tropicales y subtropicales"
$ws.Range("I469").Value = 100108005
$ws.Range("J469").Value = "Piña"
$ws.Range("K469").Value = "Caramelo"
$ws.Range("L469").Value = "Tercera"
$ws.Range("M469").Value = 210
$ws.Range("N469").Value = 23000
$ws.Range("O469").Value = 23000
$ws.Range("P469").Value = 23000
$ws.Range("Q469").Value = "$/caja 16 unidades"
$ws.Range("R469").Value = "Ecuador"
$ws.Range("S469").Value = 1438
$ws.Range("T469").Value = 16
